# Add Width/Height columns (H, I) to the "CADs" sheet, populate header +
# values, and move the active selection as shown in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CADs")

# Headers
$ws.Range("H1").Value = "Width"
$ws.Range("I1").Value = "Height"

# Data values (H2:I7)
# Row 2 is a "customFormat" row whose default style is bold (it carries the
# header-like style from row 1); the rest of the data rows (3-7) use the
# regular, non-bold body style, so make row 2's new cells match by
# explicitly turning bold back off.
$ws.Range("H2").Value = 760
$ws.Range("H2").Font.Bold = $false
$ws.Range("I2").Value = 680
$ws.Range("I2").Font.Bold = $false

$ws.Range("H3").Value = 100
$ws.Range("I3").Value = 200

$ws.Range("H4").Value = 760
$ws.Range("I4").Value = 680

$ws.Range("H5").Value = 760
$ws.Range("I5").Value = 680

$ws.Range("H6").Value = 760
$ws.Range("I6").Value = 680

$ws.Range("H7").Value = 760
$ws.Range("I7").Value = 680

# Move the selection as recorded in the saved file.
$ws.Range("K10").Select()
